# Update "想去人数" (F column) counts that changed when the site was
# regenerated, on both the "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 7645
$ws1.Range("F6").Value  = 5558
$ws1.Range("F11").Value = 246
$ws1.Range("F12").Value = 191
$ws1.Range("F13").Value = 48

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 7645
$ws4.Range("F6").Value  = 5558
$ws4.Range("F11").Value = 246
$ws4.Range("F14").Value = 191
$ws4.Range("F15").Value = 48
